$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31 (2019-12-20) ---
# Seed A31 by copying an existing date-formatted cell (A30) so it reuses the
# existing numFmtId 14 style (s="1") instead of minting a new xf.
$ws.Range("A30").Copy($ws.Range("A31")) | Out-Null
$ws.Cells.Item(31, 1).Value = 43819
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 2.5
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 1
$ws.Cells.Item(31, 12).Value = 2
$ws.Cells.Item(31, 13).Value = 3

# --- Row 32 (2019-12-21) ---
$ws.Range("A30").Copy($ws.Range("A32")) | Out-Null
$ws.Cells.Item(32, 1).Value = 43820
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 3.5
$ws.Cells.Item(32, 8).Value = 40
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 2
$ws.Cells.Item(32, 13).Value = 3

# --- Row 33 (2019-12-22) ---
# This row's date cell picks up a brand-new style (numFmtId 16, "d-mmm"),
# so set it directly via NumberFormat rather than copying an existing style.
$ws.Cells.Item(33, 1).Value = 43821
$ws.Range("A33").NumberFormat = "d-mmm"
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 15
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 35
$ws.Cells.Item(33, 6).Value = 12.5
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 1
$ws.Cells.Item(33, 12).Value = 3
$ws.Cells.Item(33, 13).Value = 3

# Match the author's final selection/active cell.
$ws.Range("N33").Select() | Out-Null
